$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-8) get re-permuted: each target row's D/J/K/L/M/O/P
# values are replaced by the values that used to live in a different source
# row. Capture the original values first, then write them back in the new
# order so we don't clobber a source row before it's been read.

$sourceForTarget = @{ 2 = 8; 3 = 4; 4 = 5; 5 = 2; 6 = 7; 7 = 6; 8 = 3 }

$original = @{}
foreach ($row in 2..8) {
    $original[$row] = @{
        D = $ws.Range("D$row").Value2
        J = $ws.Range("J$row").Value2
        K = $ws.Range("K$row").Value2
        L = $ws.Range("L$row").Value2
        M = $ws.Range("M$row").Value2
        O = $ws.Range("O$row").Value2
        P = $ws.Range("P$row").Value2
    }
}

foreach ($target in $sourceForTarget.Keys) {
    $source = $sourceForTarget[$target]
    $vals = $original[$source]
    $ws.Range("D$target").Value = $vals.D
    $ws.Range("J$target").Value = $vals.J
    $ws.Range("K$target").Value = $vals.K
    $ws.Range("L$target").Value = $vals.L
    $ws.Range("M$target").Value = $vals.M
    $ws.Range("O$target").Value = $vals.O
    $ws.Range("P$target").Value = $vals.P
}
